$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja3")

# --- Row 16: two small helper formulas ---
$ws.Range("B16").Formula = "=7*72"
$ws.Range("C16").Formula = "=523/7"

# --- Row 19: totals row ---
$ws.Range("F19").Value = 110
$ws.Range("G19").Value = 70
$ws.Range("H19").Value = 70
$ws.Range("I19").Value = 30
$ws.Range("J19").Value = 40
$ws.Range("K19").Value = 40
$ws.Range("L19").Value = 144
$ws.Range("B19").Formula = "=SUM(C19:L19)"

# --- Row 20: headers ---
$ws.Range("K20").Value = "Date"
$ws.Range("C20").Value = "CompanyName"
$ws.Range("D20").Value = "BankName"
$ws.Range("E20").Value = "AccountNumber"
$ws.Range("F20").Value = "BeneficiaryName"
$ws.Range("G20").Value = "ReportTypeName"
$ws.Range("H20").Value = "CityName"
$ws.Range("I20").Value = "ChequeNumber"
$ws.Range("J20").Value = "Amount"
$ws.Range("L20").Value = "PaymentDetail"

# --- Rows 21-25: sample cheque report rows ---
$ws.Range("C21").Value = "Ecuatoriana de Comercio S.A."
$ws.Range("D21").Value = "Banco del Pacífico"
$ws.Range("E21").Value = 11223344556677
$ws.Range("F21").Value = "DARWIN RODOLFO SANCHEZ CORREA"
$ws.Range("G21").Value = "Reporte Individual"
$ws.Range("H21").Value = " Amaluza"
$ws.Range("I21").Value = 1977
$ws.Range("J21").Value = 1453.33
$ws.Range("K21").Value = 45836
$ws.Range("L21").Value = "Detalle  general de cheques"

$ws.Range("C22").Value = "Ecuatoriana de Comercio S.A."
$ws.Range("D22").Value = "Banco del Pacífico"
$ws.Range("E22").Value = 11223344556677
$ws.Range("F22").Value = "NOMBRES Y APELLIDOS 1"
$ws.Range("G22").Value = "Reporte Individual"
$ws.Range("H22").Value = " Amaluza"
$ws.Range("I22").Value = 1978
$ws.Range("J22").Value = 123.34
$ws.Range("K22").Value = 45836
$ws.Range("L22").Value = "Detalle  general de cheques"

$ws.Range("C23").Value = "Ecuatoriana de Comercio S.A."
$ws.Range("D23").Value = "Banco del Pacífico"
$ws.Range("E23").Value = 11223344556677
$ws.Range("F23").Value = "NOMBRES Y APELLIDOS 2"
$ws.Range("G23").Value = "Reporte Individual"
$ws.Range("H23").Value = " Amaluza"
$ws.Range("I23").Value = 1979
$ws.Range("J23").Value = 223.34
$ws.Range("K23").Value = 45836
$ws.Range("L23").Value = "Detalle  general de cheques"

$ws.Range("C24").Value = "Ecuatoriana de Comercio S.A."
$ws.Range("D24").Value = "Banco del Pacífico"
$ws.Range("E24").Value = 11223344556677
$ws.Range("F24").Value = "NOMBRES Y APELLIDOS 3"
$ws.Range("G24").Value = "Reporte Individual"
$ws.Range("H24").Value = " Amaluza"
$ws.Range("I24").Value = 1980
$ws.Range("J24").Value = 323.33999999999997
$ws.Range("K24").Value = 45836
$ws.Range("L24").Value = "Detalle  general de cheques"

$ws.Range("C25").Value = "Ecuatoriana de Comercio S.A."
$ws.Range("D25").Value = "Banco del Pacífico"
$ws.Range("E25").Value = 11223344556677
$ws.Range("F25").Value = "NOMBRES Y APELLIDOS 4"
$ws.Range("G25").Value = "Reporte Individual"
$ws.Range("H25").Value = " Amaluza"
$ws.Range("I25").Value = 1981
$ws.Range("J25").Value = 423.34
$ws.Range("K25").Value = 45836
$ws.Range("L25").Value = "Detalle  general de cheques"

# Apply the date number format to K21 then fan it out (format-only paste)
# so all five cells share a single new style entry instead of five.
$ws.Range("K21").NumberFormat = "mm-dd-yy"
$ws.Range("K21").Copy()
$ws.Range("K22:K25").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- view state tweaks ---
$ws.Application.ActiveWindow.ScrollRow = 11
$ws.Range("F15").Select()
